$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 65; this pushes the existing rows 65-118
# down to 67-120, preserving all of their data (including the last two rows
# which end up as brand-new rows 119-120).
$ws.Rows("65:66").Insert()

# Populate the two freshly inserted rows with the new price-report entries.
$ws.Range("A65").Value = 12
$ws.Range("B65").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C65").Value = "Metropolitana"
$ws.Range("D65").Value = 44435
$ws.Range("E65").Value = 13
$ws.Range("F65").Value = 100112043
$ws.Range("G65").Value = "Pepino dulce"
$ws.Range("H65").Value = "Cultivar IV Región"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 300
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = 20000
$ws.Range("N65").Value = "`$/bandeja 18 kilos"
$ws.Range("O65").Value = "Provincia de Limarí"
$ws.Range("P65").Value = 1111
$ws.Range("Q65").Value = 18
$ws.Range("R65").Value = "Hortaliza"

$ws.Range("A66").Value = 12
$ws.Range("B66").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C66").Value = "Metropolitana"
$ws.Range("D66").Value = 44431
$ws.Range("E66").Value = 13
$ws.Range("F66").Value = 100112043
$ws.Range("G66").Value = "Pepino dulce"
$ws.Range("H66").Value = "Cultivar IV Región"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 300
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = 20000
$ws.Range("N66").Value = "`$/bandeja 18 kilos"
$ws.Range("O66").Value = "Provincia de Limarí"
$ws.Range("P66").Value = 1111
$ws.Range("Q66").Value = 18
$ws.Range("R66").Value = "Hortaliza"
